$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 updates: new velocity (v) -> now solve with B2=100 (time?), B3=0 ---
$ws1.Range("B2").Value = 100
$ws1.Range("B3").Value = 0
$ws1.Range("E14").Formula = "=F8"
$ws1.Range("F17").Select() | Out-Null

# --- Add Sheet2 (after Sheet1) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "g"
$ws2.Range("B1").Value = 10
$ws2.Range("A2").Value = "x"
$ws2.Range("B2").Value = 30
$ws2.Range("A3").Value = "y"
$ws2.Range("B3").Value = 0
$ws2.Range("A4").Value = "a"
$ws2.Range("B4").Value = 45
$ws2.Range("A5").Value = "v"

$ws2.Range("A13").Value = "yhit="

$ws2.Range("B10").Value = "xtan(a)"
$ws2.Range("E10").Value = "2u^2"

$ws2.Range("B5").Value = 33

$ws2.Range("C10").Formula = "=B2*TAN(B4)"
$ws2.Range("F10").Formula = "=2*POWER(B5,2)"

$ws2.Range("B11").Value = "gx^2"
$ws2.Range("C11").Formula = "=B1 * POWER(B2,2)"
$ws2.Range("E11").Value = "cos^2(a)"
$ws2.Range("F11").Formula = "=POWER(COS(B4),2)"

$ws2.Range("B13").Formula = "=C10 - (C11 / F10*F11)"

$ws2.Columns.Item(2).ColumnWidth = 9.7109375

$ws2.Range("C13").Select() | Out-Null

# --- Add Sheet3 (after Sheet2) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

$ws3.Range("A1").Value = "g"
$ws3.Range("B1").Value = 9.81
$ws3.Range("A2").Value = "x"
$ws3.Range("B2").Value = 100
$ws3.Range("A3").Value = "y"
$ws3.Range("B3").Value = 10
$ws3.Range("A4").Value = "a"
$ws3.Range("A5").Value = "v"
$ws3.Range("B5").Value = 33.014282

$ws3.Range("A9").Value = "v2"
$ws3.Range("B9").Formula = "=B5*B5"

$ws3.Range("A10").Value = "v4"
$ws3.Range("B10").Formula = "=B5*B5*B5*B5"

$ws3.Range("A11").Value = "x2"
$ws3.Range("B11").Formula = "=B2*B2"

$ws3.Range("A12").Value = "sqrt"
$ws3.Range("B12").Formula = "=SQRT(B10-B1*(B1*B11+2*B3*B9))"

$ws3.Range("A14").Formula = "=B9+B12"
$ws3.Range("B14").Formula = "=A14/ (B1*B2)"
$ws3.Range("C14").Formula = "=ATAN(B14)"
$ws3.Range("D14").Formula = "=DEGREES(C14)"

$ws3.Range("A15").Formula = "=B9-B12"
$ws3.Range("B15").Formula = "=A15/(B1*B2)"
$ws3.Range("C15").Formula = "=ATAN(B15)"
$ws3.Range("D15").Formula = "=DEGREES(C15)"

$ws3.Range("E24").Select() | Out-Null
$ws3.Activate() | Out-Null
